$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.060.23"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.650.74"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5214"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2636"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.80%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06328"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("E10").Value = "  -0.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07668"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.585"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("D13").Value = "1.640.46"
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "1.877.23"
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5597"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("D16").Value = "0.0₅8138"
$ws.Range("E16").Value = "  +2.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").Value = "26.045.98"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.621"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("E21").Value = "  +4.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "191.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.07%  "
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.004"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1186"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.49%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.521"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05438"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.269"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.446"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.343"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.555"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.50%  "
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("E36").Value = "  -0.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9445"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5632"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01580"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.858"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "1.025.78"
$ws.Range("E42").Value = "  -3.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8264"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.05%  "
$ws.Range("D45").Value = "1.786.10"
$ws.Range("E45").Value = "  -0.05%  "
$ws.Range("E46").Value = "  +5.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9984"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4331"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.958"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05138"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.37%  "
